# Update Covid-19 "paises" dashboard: refresh country stats and fix a few
# country-name/row alignments (Ghana/Azerbaiyan/Guatemala, Cabo Verde/Malaui/
# Hong Kong/Benin/Tunez, Laos/Santa Lucia, Islas Malvinas/Groenlandia) plus the
# "last updated" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last-updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 17:23"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2691508
$ws.Range("C4").Value = 9697
$ws.Range("D4").Value = 1122655
$ws.Range("E4").Value = 1439890
$ws.Range("G4").Value = 180
$ws.Range("H4").Value = 128963

# Row 5: Brasil
$ws.Range("B5").Value = 1373006
$ws.Range("C5").Value = 2518
$ws.Range("E5").Value = 557138
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 58406

# Row 7: India
$ws.Range("B7").Value = 574926
$ws.Range("C7").Value = 7390
$ws.Range("D7").Value = 340225
$ws.Range("E7").Value = 217663
$ws.Range("G7").Value = 134
$ws.Range("H7").Value = 17038

# Row 9: España
$ws.Range("B9").Value = 296351
$ws.Range("C9").Value = 301
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = 28355

# Row 28: Argentina
$ws.Range("D28").Value = 22028
$ws.Range("E28").Value = 38957
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 1283

# Row 38: Singapur
$ws.Range("D38").Value = 38500
$ws.Range("E38").Value = 5381

# Row 58: Ghana
$ws.Range("A58").Value = "Ghana"
$ws.Range("B58").Value = 17741
$ws.Range("C58").Value = 390
$ws.Range("D58").Value = 13268
$ws.Range("E58").Value = 4361
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 112

# Row 59: Azerbaiyan
$ws.Range("A59").Value = "Azerbaiyan"
$ws.Range("B59").Value = 17524
$ws.Range("C59").Value = 556
$ws.Range("D59").Value = 9715
$ws.Range("E59").Value = 7596
$ws.Range("G59").Value = 7
$ws.Range("H59").Value = 213

# Row 60: Guatemala
$ws.Range("A60").Value = "Guatemala"
$ws.Range("B60").Value = 17409
$ws.Range("C60").Value = 479
$ws.Range("D60").Value = 3170
$ws.Range("E60").Value = 13493
$ws.Range("G60").Value = 19
$ws.Range("H60").Value = 746

# Row 61: Moldavia
$ws.Range("B61").Value = 16613
$ws.Range("C61").Value = 256
$ws.Range("E61").Value = 6686
$ws.Range("G61").Value = 9
$ws.Range("H61").Value = 545

# Row 90: Republica de Yibuti
$ws.Range("B90").Value = 4682
$ws.Range("C90").Value = 26
$ws.Range("D90").Value = 4524
$ws.Range("E90").Value = 104
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 54

# Row 97: Grecia
$ws.Range("B97").Value = 3409
$ws.Range("C97").Value = 19
$ws.Range("E97").Value = 1843
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 192

# Row 124: Cabo Verde
$ws.Range("A124").Value = "Cabo Verde"
$ws.Range("B124").Value = 1226
$ws.Range("C124").Value = 61
$ws.Range("D124").Value = 608
$ws.Range("E124").Value = 604
$ws.Range("G124").Value = 2

# Row 125: Malaui
$ws.Range("A125").Value = "Malaui"
$ws.Range("B125").Value = 1224
$ws.Range("C125").Value = 72
$ws.Range("D125").Value = 260
$ws.Range("E125").Value = 950
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 14

# Row 126: Hong Kong
$ws.Range("A126").Value = "Hong Kong"
$ws.Range("B126").Value = 1206
$ws.Range("C126").Value = 2
$ws.Range("D126").Value = 1107
$ws.Range("E126").Value = 92
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 7

# Row 127: Benin
$ws.Range("A127").Value = "Benin"
$ws.Range("B127").Value = 1199
$ws.Range("C127").Value = 12
$ws.Range("D127").Value = 333
$ws.Range("E127").Value = 845
$ws.Range("G127").Value = 2
$ws.Range("H127").Value = 21

# Row 128: Tunez
$ws.Range("A128").Value = "Tunez"
$ws.Range("B128").Value = 1172
$ws.Range("D128").Value = 1029
$ws.Range("E128").Value = 93
$ws.Range("H128").Value = 50

# Row 161: Birmania
$ws.Range("D161").Value = 222
$ws.Range("E161").Value = 71

# Row 179: Trinidad yTobago
$ws.Range("B179").Value = 130
$ws.Range("C179").Value = 4
$ws.Range("D179").Value = 113

# Row 203: Laos
$ws.Range("A203").Value = "Laos"

# Row 204: Santa Lucia
$ws.Range("A204").Value = "Santa Lucia"

# Row 209: Islas Malvinas
$ws.Range("A209").Value = "Islas Malvinas"

# Row 210: Groenlandia
$ws.Range("A210").Value = "Groenlandia"
